# Datentreiber.xlsx - "Testdaten" sheet: replace the "Testmanager" row with
# "Testanalyst" in A2 and drop the now-duplicate row that used to live in A3
# (the sheet only needs one "Testanalyst" entry, not both roles).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testdaten")

# A2 held "Testmanager" -> change it to "Testanalyst"
$ws.Range("A2").Value = "Testanalyst"

# A3 held the original "Testanalyst" row - remove the whole row now that
# A2 carries that value, shrinking the used range back down to A1:A2
$ws.Rows.Item(3).Delete() | Out-Null

# leave the selection on row 2, matching the last interactive edit
$ws.Rows.Item(2).Select() | Out-Null
